$d = $word.ActiveDocument

# 1) Remove the empty paragraph right after "CHIAVÓN, Cristian:" (paraId 3CD6C6A1)
$d.Paragraphs.Item(15).Range.Delete()


# 2) Merge runs for "Sea C={...} una colección de conjuntos." paragraph (remove proofErr wrapping)
$p = $d.Paragraphs.Item(24)
$r = $p.Range
$full = $r.Text.Substring(0, $r.Text.Length - 1)
[void]$r.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, $full, 2)
